# Update the timestamp embedded in the test email addresses on the
# "UsuariosRegistro" sheet (column C: E-Mail), replacing the old
# "20251111_202811" stamp with the new "20251112_211458" stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UsuariosRegistro")

$oldStamp = "20251111_202811"
$newStamp = "20251112_211458"

# Data rows are 2..6 in column C ("E-Mail")
for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $value = $cell.Value2
    if ($value -ne $null -and $value.ToString().Contains($oldStamp)) {
        $cell.Value2 = $value.ToString().Replace($oldStamp, $newStamp)
    }
}

$wb.Save()
